$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.495.08"
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").Value = "2.242.48"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "319.94"
$ws.Range("E5").Value = "  +4.11%  "
$ws.Range("D6").Value = "101.04"
$ws.Range("E6").Value = "  +3.97%  "
$ws.Range("D7").Value = "0.585"
$ws.Range("E7").Value = "  +3.43%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").Value = "0.565"
$ws.Range("E9").Value = "  +2.88%  "
$ws.Range("D10").Value = "37.56"
$ws.Range("E10").Value = "  +3.97%  "
$ws.Range("D11").Value = "0.0836"
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("D12").Value = "7.74"
$ws.Range("E12").Value = "  +3.30%  "
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  +3.22%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.870"
$ws.Range("E14").Value = "  +2.59%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.581.64"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "14.37"
$ws.Range("E16").Value = "  +4.82%  "
$ws.Range("D17").Value = "2.245.42"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").Value = "43.437.30"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("D19").Value = "14.26"
$ws.Range("E19").Value = "  +3.57%  "
$ws.Range("D20").Value = "0.0₃0983"
$ws.Range("E20").Value = "  +5.54%  "
$ws.Range("D21").Value = "6.66"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("D22").Value = "65.71"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").Value = "3.17"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").Value = "237.45"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  +5.16%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +2.79%  "
$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  +4.31%  "
$ws.Range("D30").Value = "6.41"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").Value = "36.63"
$ws.Range("E31").Value = "  +11.79%  "
$ws.Range("D32").Value = "20.35"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").Value = "0.0876"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "160.49"
$ws.Range("E34").Value = "  +3.29%  "
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("D36").Value = "3.24"
$ws.Range("E36").Value = "  +3.57%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  +6.43%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "0.105"
$ws.Range("E40").Value = "  +1.30%  "
$ws.Range("D41").Value = "3.74"
$ws.Range("E41").Value = "  +8.66%  "
$ws.Range("D42").Value = "0.0324"
$ws.Range("E42").Value = "  +3.16%  "
$ws.Range("D43").Value = "14.73"
$ws.Range("E43").Value = "  +26.59%  "
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "1.827.72"
$ws.Range("E45").Value = "  +3.61%  "
$ws.Range("D46").Value = "0.206"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("E47").Value = "  -3.14%  "
$ws.Range("D48").Value = "5.32"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "75.10"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "8.79"
$ws.Range("E50").Value = "  +4.51%  "
$ws.Range("D51").Value = "58.98"
$ws.Range("E51").Value = "  -0.25%  "
